$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Пройдено" (actually traveled) text for row 16 to include the new summit
$ws.Range("D16").Value = "пер. Col de Nannaz (1А, 2773) – вер. Becca Trecare (н/к,  3032) – Col de Fontaines (н/к, 2696) -- Mont-Perron -- camping Glair"

# Update the height delta (Δh, м) for row 16
$ws.Range("H16").Value = "+450,-1600"

# Update row 16 height to fit the new wrapped text
$ws.Rows.Item(16).RowHeight = 124.6

# Move the visible top-left cell and active selection like the author left it
$ws.Range("D16").Select()
$excel.ActiveWindow.ScrollRow = 13
